$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 685.3782603333333
$ws.Range("H2").Value = 2056.134781
$ws.Range("I2").Value = 0.7131285654702259
$ws.Range("J2").Value = 0.7131285654702259
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.288321666666667
$ws.Range("N2").Value = 3.864965
$ws.Range("Q2").Value = 882.9876626497406
$ws.Range("R2").Value = 7946.888963847665
$ws.Range("S2").Value = 0.7131285654702259
$ws.Range("T2").Value = 0.7131285654702259

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 213.8079683333333
$ws.Range("H3").Value = 641.423905
$ws.Range("I3").Value = 0.2224648468854243
$ws.Range("J3").Value = 0.2224648468854243
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.288321666666667
$ws.Range("N3").Value = 3.864965
$ws.Range("Q3").Value = 275.4534381098139
$ws.Range("R3").Value = 2479.080942988325
$ws.Range("S3").Value = 0.2224648468854243
$ws.Range("T3").Value = 0.2224648468854243

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 61.90030400000001
$ws.Range("H4").Value = 185.700912
$ws.Range("I4").Value = 0.06440658764434989
$ws.Range("J4").Value = 0.06440658764434989
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.288321666666667
$ws.Range("N4").Value = 3.864965
$ws.Range("Q4").Value = 79.74750281645335
$ws.Range("R4").Value = 717.7275253480801
$ws.Range("S4").Value = 0.06440658764434989
$ws.Range("T4").Value = 0.06440658764434989
